$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gains two new data rows ("line7", "line8") which sit logically
# between "line6" and "extr1" (row 7 and the old row 8), pushing extr1..extr8
# down by two rows (from rows 8-15 to rows 10-17). Shift the existing extr
# rows down first, working from the bottom up so source data isn't
# clobbered. Each shift pastes values then formats so the new destination
# row picks up the bold/border/centered "A column" style without minting a
# new style entry (PasteSpecial(xlPasteAll) does not propagate formats to a
# brand new row, so values + formats are pasted separately).
function Move-Row($src, $dst) {
    $srcRange = "A" + $src + ":E" + $src
    $dstCell = "A" + $dst
    $ws.Range($srcRange).Copy() | Out-Null
    $ws.Range($dstCell).PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $ws.Range($srcRange).Copy() | Out-Null
    $ws.Range($dstCell).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

Move-Row 15 17
Move-Row 14 16
Move-Row 13 15
Move-Row 12 14
Move-Row 11 13
Move-Row 10 12
Move-Row 9 11
Move-Row 8 10
$excel.CutCopyMode = $false

# Rows 10-17 now hold extr1..extr8 (shifted copies of the old rows 8-15).
# Write the two brand new rows (line7, line8) into rows 8 and 9, copying
# formatting from the row above (line6, row 7) so column A keeps its
# bold/border/centered style.
$ws.Range("A7:E7").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Renumber column A for the shifted extr rows and fix the in_service flags
# that changed in the source data.
$ws.Range("A10").Value = 8
$ws.Range("E10").Value = $true

$ws.Range("A11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("A12").Value = 10
$ws.Range("E12").Value = $true

$ws.Range("A13").Value = 11

$ws.Range("A14").Value = 12
$ws.Range("E14").Value = $false

$ws.Range("A15").Value = 13
$ws.Range("E15").Value = $false

$ws.Range("A16").Value = 14

$ws.Range("A17").Value = 15
$ws.Range("E17").Value = $false
